{"js": "// Update the division-problem answers in the single table of the\n// worksheet. Every populated cell (\"A\u00f7B=C, D\") is replaced with a new\n// value per the mapping below (order = reading order: row by row,\n// left to right). Blank spacer rows/cells are left untouched.\nconst replacements = new Map([\n  [\"768\u00f75=153, 3\", \"273\u00f77=39, 0\"],\n  [\"690\u00f79=76, 6\", \"144\u00f72=72, 0\"],\n  [\"350\u00f79=38, 8\", \"653\u00f75=130, 3\"],\n  [\"240\u00f77=34, 2\", \"554\u00f74=138, 2\"],\n  [\"355\u00f77=50, 5\", \"530\u00f77=75, 5\"],\n  [\"527\u00f72=263, 1\", \"620\u00f77=88, 4\"],\n  [\"268\u00f75=53, 3\", \"790\u00f74=197, 2\"],\n  [\"274\u00f77=39, 1\", \"100\u00f72=50, 0\"],\n  [\"777\u00f73=259, 0\", \"414\u00f75=82, 4\"],\n  [\"116\u00f72=58, 0\", \"639\u00f76=106, 3\"],\n  [\"874\u00f78=109, 2\", \"374\u00f72=187, 0\"],\n  [\"188\u00f75=37, 3\", \"232\u00f75=46, 2\"],\n  [\"438\u00f74=109, 2\", \"123\u00f72=61, 1\"],\n  [\"793\u00f79=88, 1\", \"270\u00f74=67, 2\"],\n  [\"660\u00f76=110, 0\", \"115\u00f73=38, 1\"],\n  [\"850\u00f78=106, 2\", \"749\u00f79=83, 2\"],\n  [\"821\u00f74=205, 1\", \"847\u00f73=282, 1\"],\n  [\"968\u00f74=242, 0\", \"791\u00f73=263, 2\"],\n  [\"880\u00f75=176, 0\", \"399\u00f78=49, 7\"],\n  [\"659\u00f78=82, 3\", \"937\u00f77=133, 6\"],\n  [\"589\u00f77=84, 1\", \"140\u00f72=70, 0\"],\n  [\"195\u00f77=27, 6\", \"462\u00f76=77, 0\"],\n  [\"956\u00f76=159, 2\", \"971\u00f76=161, 5\"],\n  [\"843\u00f78=105, 3\", \"776\u00f79=86, 2\"],\n  [\"848\u00f72=424, 0\", \"735\u00f73=245, 0\"],\n]);\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  table.load(\"values\");\n  await context.sync();\n\n  const values = table.values;\n  let changed = false;\n\n  for (let r = 0; r < values.length; r++) {\n    for (let c = 0; c < values[r].length; c++) {\n      const cellText = values[r][c];\n      if (replacements.has(cellText)) {\n        values[r][c] = replacements.get(cellText);\n        changed = true;\n      }\n    }\n  }\n\n  if (changed) {\n    table.values = values;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem answers in the single table of the\n# worksheet. Every populated cell (\"A\u00f7B=C, D\") is replaced with a new\n# value per the mapping below. Blank spacer rows/cells are untouched.\n# Uses Find/Replace over the whole document range so cell formatting\n# (fonts/size/alignment) carried on the run is left intact.\n\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n    \"768\u00f75=153, 3\" = \"273\u00f77=39, 0\"\n    \"690\u00f79=76, 6\"  = \"144\u00f72=72, 0\"\n    \"350\u00f79=38, 8\"  = \"653\u00f75=130, 3\"\n    \"240\u00f77=34, 2\"  = \"554\u00f74=138, 2\"\n    \"355\u00f77=50, 5\"  = \"530\u00f77=75, 5\"\n    \"527\u00f72=263, 1\" = \"620\u00f77=88, 4\"\n    \"268\u00f75=53, 3\"  = \"790\u00f74=197, 2\"\n    \"274\u00f77=39, 1\"  = \"100\u00f72=50, 0\"\n    \"777\u00f73=259, 0\" = \"414\u00f75=82, 4\"\n    \"116\u00f72=58, 0\"  = \"639\u00f76=106, 3\"\n    \"874\u00f78=109, 2\" = \"374\u00f72=187, 0\"\n    \"188\u00f75=37, 3\"  = \"232\u00f75=46, 2\"\n    \"438\u00f74=109, 2\" = \"123\u00f72=61, 1\"\n    \"793\u00f79=88, 1\"  = \"270\u00f74=67, 2\"\n    \"660\u00f76=110, 0\" = \"115\u00f73=38, 1\"\n    \"850\u00f78=106, 2\" = \"749\u00f79=83, 2\"\n    \"821\u00f74=205, 1\" = \"847\u00f73=282, 1\"\n    \"968\u00f74=242, 0\" = \"791\u00f73=263, 2\"\n    \"880\u00f75=176, 0\" = \"399\u00f78=49, 7\"\n    \"659\u00f78=82, 3\"  = \"937\u00f77=133, 6\"\n    \"589\u00f77=84, 1\"  = \"140\u00f72=70, 0\"\n    \"195\u00f77=27, 6\"  = \"462\u00f76=77, 0\"\n    \"956\u00f76=159, 2\" = \"971\u00f76=161, 5\"\n    \"843\u00f78=105, 3\" = \"776\u00f79=86, 2\"\n    \"848\u00f72=424, 0\" = \"735\u00f73=245, 0\"\n}\n\nforeach ($key in $replacements.Keys) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $key\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replacements[$key]\n    $find.Execute($key, $false, $false, $false, $false, $false, $true, 1, $false, $replacements[$key], 2)\n}\n"}
